$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44281
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 31000
$ws.Range("L2").Value = 32000
$ws.Range("M2").Value = 31500
$ws.Range("P2").Value = 1260

$ws.Range("D3").Value = 44267
$ws.Range("J3").Value = 45
$ws.Range("M3").Value = 24333
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 973

$ws.Range("D4").Value = 44259
$ws.Range("J4").Value = 65
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24538
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 982

$ws.Range("D5").Value = 44181
$ws.Range("K5").Value = 19500
$ws.Range("M5").Value = 19750
$ws.Range("P5").Value = 790

$ws.Range("D6").Value = 44249
$ws.Range("K6").Value = 21000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21500
$ws.Range("P6").Value = 860

$ws.Range("D7").Value = 44176
$ws.Range("K7").Value = 19000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 19500
$ws.Range("P7").Value = 780

$ws.Range("D8").Value = 44208
$ws.Range("J8").Value = 65
$ws.Range("K8").Value = 22000
$ws.Range("M8").Value = 23385
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 935

$ws.Range("D9").Value = 44179
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 16000
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 16500
$ws.Range("N9").Value = "`$/saco 25 kilos"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 660

$ws.Range("D11").Value = 44323
$ws.Range("K11").Value = 29000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 29500
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 1180

$ws.Range("D12").Value = 44160
$ws.Range("H12").Value = "Magnum"
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 30000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 30000
$ws.Range("N12").Value = "`$/malla 25 kilos"
$ws.Range("O12").Value = "Región de Coquimbo"
$ws.Range("P12").Value = 1200

$ws.Range("D13").Value = 44160
$ws.Range("H13").Value = "Magnum"
$ws.Range("J13").Value = 35
$ws.Range("K13").Value = 28000
$ws.Range("L13").Value = 28000
$ws.Range("M13").Value = 28000
$ws.Range("P13").Value = 1120

$ws.Range("D14").Value = 44211
$ws.Range("J14").Value = 70
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 23143
$ws.Range("P14").Value = 926

$ws.Range("D15").Value = 44272
$ws.Range("J15").Value = 42
$ws.Range("L15").Value = 24000
$ws.Range("M15").Value = 22857
$ws.Range("O15").Value = "Región de O'Higgins"
$ws.Range("P15").Value = 914

$ws.Range("D16").Value = 44218
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24562
$ws.Range("P16").Value = 982

$ws.Range("D17").Value = 44210
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = 23000
$ws.Range("M17").Value = 23857
$ws.Range("P17").Value = 954

$ws.Range("D19").Value = 44174
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = 21000
$ws.Range("L19").Value = 22000
$ws.Range("M19").Value = 21500
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 860

$ws.Range("D20").Value = 44244
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 25000
$ws.Range("L20").Value = 26000
$ws.Range("M20").Value = 25500
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 1020

$ws.Range("D21").Value = 44195
$ws.Range("J21").Value = 160
$ws.Range("K21").Value = 32000
$ws.Range("L21").Value = 33000
$ws.Range("M21").Value = 32500
$ws.Range("P21").Value = 1300

$ws.Range("D22").Value = 44209
$ws.Range("J22").Value = 90
$ws.Range("M22").Value = 23889
$ws.Range("P22").Value = 956

$ws.Range("D23").Value = 44232
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = 24000
$ws.Range("L23").Value = 25000
$ws.Range("M23").Value = 24500
$ws.Range("P23").Value = 980

$ws.Range("D24").Value = 44273
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 24000
$ws.Range("M24").Value = 23455
$ws.Range("P24").Value = 938

$ws.Range("D25").Value = 44236
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 22000
$ws.Range("L25").Value = 23000
$ws.Range("M25").Value = 22500
$ws.Range("O25").Value = "Provincia de Diguillín"
$ws.Range("P25").Value = 900

$ws.Range("D26").Value = 44203
$ws.Range("J26").Value = 50
$ws.Range("L26").Value = 24000
$ws.Range("M26").Value = 23200
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 928

$ws.Range("D27").Value = 44168
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 15000
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 15500
$ws.Range("N27").Value = "`$/saco 25 kilos"
$ws.Range("O27").Value = "Región del Maule"
$ws.Range("P27").Value = 620

$ws.Range("D28").Value = 44161
$ws.Range("H28").Value = "Magnum"
$ws.Range("J28").Value = 47
$ws.Range("K28").Value = 28000
$ws.Range("L28").Value = 29000
$ws.Range("M28").Value = 28532
$ws.Range("O28").Value = "Región de O'Higgins"
$ws.Range("P28").Value = 1141

$ws.Range("D29").Value = 44258
$ws.Range("J29").Value = 55
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 23909
$ws.Range("O29").Value = "Provincia de Diguillín"
$ws.Range("P29").Value = 956

$ws.Range("D30").Value = 44252
$ws.Range("J30").Value = 80
$ws.Range("K30").Value = 22000
$ws.Range("L30").Value = 23000
$ws.Range("M30").Value = 22500
$ws.Range("O30").Value = "Provincia de Diguillín"
$ws.Range("P30").Value = 900

$ws.Range("D31").Value = 44186
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("J31").Value = 60
$ws.Range("K31").Value = 19000
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = 19500
$ws.Range("N31").Value = "`$/saco 25 kilos"
$ws.Range("O31").Value = "Región del Maule"
$ws.Range("P31").Value = 780

$ws.Range("D32").Value = 44159
$ws.Range("J32").Value = 47
$ws.Range("K32").Value = 27000
$ws.Range("M32").Value = 27532
$ws.Range("N32").Value = "`$/malla 25 kilos"
$ws.Range("P32").Value = 1101

$ws.Range("D33").Value = 44166
$ws.Range("H33").Value = "Magnum"
$ws.Range("J33").Value = 38
$ws.Range("K33").Value = 24000
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = 24526
$ws.Range("N33").Value = "`$/malla 25 kilos"
$ws.Range("O33").Value = "Región de Coquimbo"
$ws.Range("P33").Value = 981

$ws.Range("D34").Value = 44334
$ws.Range("J34").Value = 30
$ws.Range("K34").Value = 31000
$ws.Range("L34").Value = 32000
$ws.Range("M34").Value = 31500
$ws.Range("O34").Value = "Región del Maule"
$ws.Range("P34").Value = 1260

$ws.Range("D35").Value = 44250
$ws.Range("J35").Value = 120
$ws.Range("K35").Value = 22000
$ws.Range("L35").Value = 23000
$ws.Range("M35").Value = 22500
$ws.Range("O35").Value = "Provincia de Diguillín"
$ws.Range("P35").Value = 900

$ws.Range("D36").Value = 44201
$ws.Range("J36").Value = 33
$ws.Range("K36").Value = 26000
$ws.Range("L36").Value = 28000
$ws.Range("M36").Value = 27091
$ws.Range("P36").Value = 1084

$ws.Range("D37").Value = 44193
$ws.Range("K37").Value = 35000
$ws.Range("L37").Value = 36000
$ws.Range("M37").Value = 35500
$ws.Range("P37").Value = 1420
